$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Status cells (D2:D4) to FAILED and change their fill to red
$ws.Range("D2").Value = "FAILED"
$ws.Range("D3").Value = "FAILED"
$ws.Range("D4").Value = "FAILED"

$ws.Range("D2").Interior.ColorIndex = 3
$ws.Range("D3").Interior.ColorIndex = 3
$ws.Range("D4").Interior.ColorIndex = 3

# Update Execution Time cells (E2:E4)
$ws.Range("E2").Value = "2096 ms"
$ws.Range("E3").Value = "3297 ms"
$ws.Range("E4").Value = "3374 ms"

# Update Failure Reason cells (F2:F4)
$ws.Range("F2").Value = "User with email containing 'tuong.2274802010979' should be found expected [true] but found [false]"
$ws.Range("F3").Value = "User with ID 2274802010979 should be found expected [true] but found [false]"
$ws.Range("F4").Value = "User with name 'Bui Ke Ton Tuong' should be found expected [true] but found [false]"

# Adjust column widths (closest values reachable via the ColumnWidth property,
# which Excel quantizes to 1/6-character increments)
$ws.Columns.Item(4).ColumnWidth = 5.7
$ws.Columns.Item(6).ColumnWidth = 81.2
